$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Icam1", "Itgax", "M2",  3, 1, 14.452709,           43.358127,           0.1476906377370901,  0.1476906377370901,  3, 1,               45.39437433333333, 136.183123, 0.9965580896473674, 0.9965580896473674, 656.0716824767355,  5904.64514229062,   0.1471822998020759,   0.1471822998020759),
    @("ECs",  "Icam1", "Itgax", "sCs", 3, 1, 14.452709,           43.358127,           0.1476906377370901,  0.1476906377370901,  1, 0.3333333333333333, 0.156783,          0.470349,   0.00344191035263268, 0.00344191035263268, 2.265939075147,     20.393451676323,    0.0005083379350142132, 0.0005083379350142133),
    @("FAPs", "Icam1", "Itgax", "M2",  3, 1, 45.91529066666667,   137.745872,          0.4692032864180593,  0.4692032864180593,  3, 1,               45.39437433333333, 136.183123, 0.9965580896473674, 0.9965580896473674, 2084.295892146473,  18758.66302931826,  0.4675883307690477,   0.4675883307690477),
    @("FAPs", "Icam1", "Itgax", "sCs", 3, 1, 45.91529066666667,   137.745872,          0.4692032864180593,  0.4692032864180593,  1, 0.3333333333333333, 0.156783,          0.470349,   0.00344191035263268, 0.00344191035263268, 7.198737016592001,  64.78863314932801,  0.001614955649011595, 0.001614955649011595),
    @("M2",   "Icam1", "Itgax", "M2",  3, 1, 36.015061,           108.045183,          0.3680339324088102,  0.3680339324088103,  3, 1,               45.39437433333333, 136.183123, 0.9965580896473674, 0.9965580896473674, 1634.881160671834,  14713.93044604651,  0.3667671926067322,   0.3667671926067323),
    @("M2",   "Icam1", "Itgax", "sCs", 3, 1, 36.015061,           108.045183,          0.3680339324088102,  0.3680339324088103,  1, 0.3333333333333333, 0.156783,          0.470349,   0.00344191035263268, 0.00344191035263268, 5.646549308763,     50.81894377886701,  0.001266739802078,    0.001266739802078),
    @("sCs",  "Icam1", "Itgax", "M2",  3, 1, 1.474929666666667,   4.424789000000001,  0.01507214343604052, 0.01507214343604052, 3, 1,               45.39437433333333, 136.183123, 0.9965580896473674, 0.9965580896473674, 66.95350940400523,  602.5815846360471,  0.01502026646951165,  0.01502026646951165),
    @("sCs",  "Icam1", "Itgax", "sCs", 3, 1, 1.474929666666667,   4.424789000000001,  0.01507214343604052, 0.01507214343604052, 1, 0.3333333333333333, 0.156783,          0.470349,   0.00344191035263268, 0.00344191035263268, 0.231243897929,     2.081195081361,     0.00005187696652887256, 0.00005187696652887256)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = 1 + $c
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}
